$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as TEXT (so it is stored as a shared string, matching
# the workbook's existing convention of storing these "numeric looking" labels as
# text), while preserving whatever NumberFormat the cell already has.
function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $origNF = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $origNF
}

# ------------------------------------------------------------------
# 1) Insert the 3 new rows in their correct (alphabetically-sorted)
#    positions. Inserting a row in this engine inherits the visual
#    formatting (number format, font, fill, alignment) of the row
#    above automatically, but NOT the row height or the merged-cell
#    state, so those are fixed up explicitly right after.
# ------------------------------------------------------------------

# Row 10 -> CLARITINE (between BETADERM at row 9 and DOSTINEX which was row 10)
$ws.Rows(10).Insert()
$ws.Rows(10).RowHeight = 24.75
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# Row 13 -> MAXILASE (after ERASTAPEX, which is now row 12)
$ws.Rows(13).Insert()
$ws.Rows(13).RowHeight = 24.75
$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()

# Row 14 -> OPLEX-N (after MAXILASE)
$ws.Rows(14).Insert()
$ws.Rows(14).RowHeight = 25.5
$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

# ------------------------------------------------------------------
# 2) Rewrite every data row (7-22) with its final target content.
#    This is simplest & safest: rows that didn't move keep the same
#    values: rows that shifted down (because of the 3 inserts above)
#    just get their correct sequential "#" (column A) and unchanged
#    item data re-asserted in their new location.
# ------------------------------------------------------------------

$rows = @(
    @{ r = 7;  a = 1;  c = "AUGMENTIN 1 GM 14 F.C. TABS.";             h = "2:1"; l = "1"; n = "210.00"; p = "105.0000"; q = "0:1" }
    @{ r = 8;  a = 2;  c = "BECOZYME I.M./I.V. 12 AMP";                h = "2:5"; l = "1"; n = "120.00"; p = "39.6000";  q = "0:4" }
    @{ r = 9;  a = 3;  c = "BETADERM 0.1% CREAM 15 GM";                h = "6:0"; l = "1"; n = "18.00";  p = "18.0000";  q = "1:0" }
    @{ r = 10; a = 4;  c = "CLARITINE 1MG/ML SYRUP 100ML";             h = "1:0"; l = "1"; n = "62.00";  p = "62.0000";  q = "1:0" }
    @{ r = 11; a = 5;  c = "DOSTINEX 0.5 MG 2 TABS.";                  h = "1:0"; l = "1"; n = "172.00"; p = "172.0000"; q = "1:0" }
    @{ r = 12; a = 6;  c = "ERASTAPEX TRIO 10/40/25MG 30 F.C. TABS.";  h = "0:1"; l = "1"; n = "162.00"; p = "162.0000"; q = "1:0" }
    @{ r = 13; a = 7;  c = "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML";    h = "3:0"; l = "1"; n = "57.00";  p = "57.0000";  q = "1:0" }
    @{ r = 14; a = 8;  c = "OPLEX-N SYRUP 125ML";                      h = "4:0"; l = "1"; n = "31.00";  p = "31.0000";  q = "1:0" }
    @{ r = 15; a = 9;  c = "PANTAZOL 40MG 14 ENTERIC COATED TAB.";     h = "1:0"; l = "1"; n = "104.00"; p = "52.0000";  q = "0:1" }
    @{ r = 16; a = 10; c = "TRIACTIN 4MG 20 TAB";                      h = "0:1"; l = "1"; n = "46.00";  p = "23.0000";  q = "0:1" }
    @{ r = 17; a = 11; c = "VIOTIC EAR DROPS 10 ML";                   h = "1:0"; l = "1"; n = "23.00";  p = "23.0000";  q = "1:0" }
    @{ r = 18; a = 12; c = "VOLTAREN 75MG/3ML 3 AMP.";                 h = "5:2"; l = "1"; n = "51.00";  p = "16.8300";  q = "0:1" }
    @{ r = 19; a = 13; c = "XORAXON 1GM I.M. VIAL";                    h = "6:0"; l = "1"; n = "56.00";  p = "112.0000"; q = "2:0" }
    @{ r = 20; a = 14; c = "ترمومتر ديجتال";                           h = "6:0"; l = "0"; n = "50.00";  p = "50.0000";  q = "1:0" }
    @{ r = 21; a = 15; c = "سرنجات 3 سم";                              h = "0:0"; l = "0"; n = "2.00";   p = "6.0000";   q = "3:0" }
    @{ r = 22; a = 16; c = "قطن 50جم";                                 h = "3:0"; l = "0"; n = "10.00";  p = "10.0000";  q = "1:0" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.a
    Set-TextValue "C$r" $row.c
    Set-TextValue "H$r" $row.h
    Set-TextValue "L$r" $row.l
    Set-TextValue "N$r" $row.n
    Set-TextValue "P$r" $row.p
    Set-TextValue "Q$r" $row.q
}

# ------------------------------------------------------------------
# 3) Totals row (was row 20, now row 23 after the 3 inserts) and the
#    footer row (was row 21, now row 24) shifted down automatically;
#    just refresh the grand-total value.
# ------------------------------------------------------------------
$ws.Range("P23").Value = 939.42999999999995

# ------------------------------------------------------------------
# 4) Update the generated-at timestamp in the footer row.
# ------------------------------------------------------------------
Set-TextValue "A24" "Monday, 1 September, 2025 10:54 AM"
